# Update the lattice-multiplication exercise table with a freshly
# generated set of problems, preserving the existing formatting
# (font size, line-break layout) of every cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry is the new content for one table cell, in row-major order
# (row 1 col 1, row 1 col 2, row 1 col 3, row 2 col 1, ...).
# Within a cell, lines are joined with a vertical-tab character, which
# Word's Range.Text setter turns into manual line breaks (<w:br/>).
$vt = [char]11

$newCells = @(
    @("21 x 14", "  1    4", "  ----", "2|    |", "1|    |"),
    @("83 x 74", "  7    4", "  ----", "8|    |", "3|    |"),
    @("23 x 65", "  6    5", "  ----", "2|    |", "3|    |"),

    @("99 x 47", "  4    7", "  ----", "9|    |", "9|    |"),
    @("64 x 27", "  2    7", "  ----", "6|    |", "4|    |"),
    @("52 x 54", "  5    4", "  ----", "5|    |", "2|    |"),

    @("60 x 78", "  7    8", "  ----", "6|    |", "0|    |"),
    @("47 x 38", "  3    8", "  ----", "4|    |", "7|    |"),
    @("86 x 22", "  2    2", "  ----", "8|    |", "6|    |"),

    @("35 x 94", "  9    4", "  ----", "3|    |", "5|    |"),
    @("13 x 30", "  3    0", "  ----", "1|    |", "3|    |"),
    @("20 x 92", "  9    2", "  ----", "2|    |", "0|    |"),

    @("15 x 93", "  9    3", "  ----", "1|    |", "5|    |"),
    @("63 x 10", "  1    0", "  ----", "6|    |", "3|    |"),
    @("98 x 19", "  1    9", "  ----", "9|    |", "8|    |")
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count

$index = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $lines = $newCells[$index]
        $text = [string]::Join($vt, $lines)
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $text
        $index = $index + 1
    }
}

Write-Output "Updated $index cells"
